{"js": "// Append 4 empty paragraphs and a new paragraph with text \"Editing file.\"\n// to the end of the document body (after the last existing paragraph,\n// which ends with \"I rise.\").\nconst body = context.document.body;\n\nbody.insertParagraph(\"\", \"End\");\nbody.insertParagraph(\"\", \"End\");\nbody.insertParagraph(\"\", \"End\");\nbody.insertParagraph(\"\", \"End\");\nbody.insertParagraph(\"Editing file.\", \"End\");\n\nawait context.sync();\n", "ps1": "# Append 4 empty paragraphs and a new paragraph with text \"Editing file.\"\n# to the end of the document (after the last existing paragraph, which\n# ends with \"I rise.\").\n$d = $word.ActiveDocument\n\n$d.Paragraphs.Add() | Out-Null\n$d.Paragraphs.Add() | Out-Null\n$d.Paragraphs.Add() | Out-Null\n$d.Paragraphs.Add() | Out-Null\n\n$p = $d.Paragraphs.Add()\n$p.Range.Text = \"Editing file.\"\n"}
